$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in columns F and G for rows 8-13 (new Hughes Fix results)
$ws.Range("F8").Value = 1.25
$ws.Range("G8").Value = 1.25

$ws.Range("F9").Value = 1.1499999999999999
$ws.Range("G9").Value = 1.05

$ws.Range("F10").Value = 1

$ws.Range("F11").Value = 0.9
$ws.Range("G11").Value = 0.8

$ws.Range("F12").Value = 0.75
$ws.Range("G12").Value = 0.75

$ws.Range("F13").Value = 0.65
$ws.Range("G13").Value = 0.55000000000000004

# Update view: scroll back to top-left A1, zoom to 115%, and move selection to I9
$ws.Activate()
$excel.ActiveWindow.Zoom = 115
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I9").Select()
